$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Acoes_compartilhadas"
$ws.Range("F1").Value = "Membros_que_executam"
$ws.Range("H1").Value = "Projetos_de_impacto"
$ws.Range("G1").Value = "Participacao_em_eventos"

$ws.Range("F5").Select()
